# evaluation excluding exact matches
#
# 1) Add a new "Sheet2" (after Sheet1) holding a second evaluation table
#    ("-exact" gold-pair metrics), mirroring Sheet1's layout but rendered
#    with an explicit black 12pt font instead of the theme font.
# 2) Resize Sheet1's column A and bump the zoom on both sheets to 200%.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- New Sheet2, appended after Sheet1 -------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "Sheet2"

# Header row (same column headers as Sheet1, reusing the same shared strings)
$ws2.Range("B1").Value = "big-small comm-comm"
$ws2.Range("C1").Value = "small-big comm-comm"
$ws2.Range("D1").Value = "big-small node-comm"
$ws2.Range("E1").Value = "small-big node-comm"

# Row labels (new shared strings)
$ws2.Range("A2").Value = "Gold pairs (-exact) in matched top 1"
$ws2.Range("A3").Value = "Gold pairs (-exact) in matched top 5"
$ws2.Range("A4").Value = "Gold pairs (-exact) in matched top 10"

# Data values
$ws2.Range("B2").Value = 0.18243243243243201
$ws2.Range("C2").Value = 0.18918918918918901
$ws2.Range("D2").Value = 0.195945945945945
$ws2.Range("E2").Value = 0.18918918918918901

$ws2.Range("B3").Value = 0.22972972972972899
$ws2.Range("C3").Value = 0.23648648648648599
$ws2.Range("D3").Value = 0.37837837837837801
$ws2.Range("E3").Value = 0.35135135135135098

$ws2.Range("B4").Value = 0.27702702702702697
$ws2.Range("C4").Value = 0.28378378378378299
$ws2.Range("D4").Value = 0.46621621621621601
$ws2.Range("E4").Value = 0.445945945945946

# Number format for the data block
$ws2.Range("B2:E4").NumberFormat = "0.00000%"

# Sheet2 uses an explicit black 12pt Calibri font on every used cell
$ws2.Range("A1:E4").Font.Size = 12
$ws2.Range("A1:E4").Font.Color = 0

# Column widths matching Sheet1's layout
$ws2.Columns.Item(1).ColumnWidth = 31.5
$ws2.Columns.Item(2).ColumnWidth = 19.666666666666668
$ws2.Columns.Item(3).ColumnWidth = 19.666666666666668
$ws2.Columns.Item(4).ColumnWidth = 18.5
$ws2.Columns.Item(5).ColumnWidth = 18.5

# --- Sheet1 tweaks -----------------------------------------------------
$ws1.Columns.Item(1).ColumnWidth = 31.5

# --- View: zoom both sheets to 200% and leave Sheet1 as the active tab -----
$ws2.Application.ActiveWindow.Zoom = 200
$ws1.Activate()
$excel.ActiveWindow.Zoom = 200

Write-Output "done"
